$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert a new row at position 78 (shifts rows 78..143 down to 79..144)
$ws.Rows.Item(78).Insert()

# 2. Grow Table1 so it covers the newly inserted row (A8:K143 -> A8:K144)
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K144"))

# 3. The inserted row (78) comes back with generic default formatting;
#    copy the correct "interior data row" formatting from row 79 (which
#    now carries what used to be row 78's formatting).
$ws.Range("A79:K79").Copy()
$ws.Range("A78:K78").PasteSpecial(-4122)

# 4. Restore the calculated-column formulas for the EARNED (helper) column
#    on the new row 78 and on the new final row 144 (row insertion can
#    leave these without a formula / with a stale one).
$ws.Cells.Item(78,7).Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Cells.Item(144,7).Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# 5. K77/K78 (REMARKS) now hold dates, so give them the same date-number
#    formatting already used by K76 above them.
$ws.Range("K76").Copy()
$ws.Range("K77:K78").PasteSpecial(-4122)

# 6. Fill in the new leave-card entries.
#    Row 77: PARTICULARS = SP(1-0-0); REMARKS = 7/21/2023
$ws.Cells.Item(77,2).Value = "SP(1-0-0)"
$ws.Cells.Item(77,11).Value = 45128

#    Row 78 (brand-new row): PARTICULARS = SL(2-0-0); Absence Undertime W/Pay = 2;
#    REMARKS = "7/6,11/2023"
$ws.Cells.Item(78,2).Value = "SL(2-0-0)"
$ws.Cells.Item(78,8).Value = 2
$ws.Cells.Item(78,11).Value = "7/6,11/2023"

# 7. Recalculate so the dependent BALANCE/summary formulas (Sheet1!I9,
#    CONVERTION!A7, etc.) pick up the new Absence Undertime value.
$excel.CalculateFullRebuild()
